# authentication-type.xlsx update: expand the 1-level Codice/Label table into a
# bilingual (IT/EN/DE) two-level (level1 + level2) lookup table, A1:H10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write all new cell text (values only - formatting is applied afterwards
#    by copying format from cells that already carry the right style, so the
#    style table doesn't grow with duplicate/near-duplicate xf entries).
# ---------------------------------------------------------------------------
$data = @(
  @('codice_1_livello', 'label_ITA_1_livello', 'label_ENG_1_livello', 'label_DEU_1_livello', 'codice_2_livello', 'label_ITA_2_livello', 'label_ENG_2_livello', 'label_DEU_2_livello'),
  @('NONE', 'Nessuna - accesso libero', 'None - free access', 'Keine - freier Zugang', $null, $null, $null, $null),
  @('SFA', 'Autenticazione singolo fattore', 'Single Factor Authentication', 'Single-Faktor-Authentifizierung', 'IDPWD', 'User ID e Password', 'User ID and password', 'Benutzer-ID und Passwort'),
  @('SFA', 'Autenticazione singolo fattore', 'Single Factor Authentication', 'Single-Faktor-Authentifizierung', 'SPIDL1', 'Credenziale SPID Livello 1', 'SPID (National public system of e-ID) Credential Level 1', 'SPID (Nationales öffentliches System der e-ID) Berechtigungsnachweis Stufe 1'),
  @('2FA', 'Autenticazione a due fattori', 'Two Factor Authentication', 'Zwei-Faktor-Authentifizierung', 'ID OTP', 'User ID e One Time Password', 'User ID and One Time Password', 'Benutzer-ID und Einmal-Passwort'),
  @('2FA', 'Autenticazione a due fattori', 'Two Factor Authentication', 'Zwei-Faktor-Authentifizierung', 'SPIDL2', 'SPID Livello 2', 'SPID (National public system of e-ID) Level 2', 'SPID (Nationales öffentliches System der e-ID) Stufe 2'),
  @('2FAHW', 'Autenticazione a due fattori hardware', 'Two Factor Authentication hardware', 'Zwei-Faktor-Hardware-Authentifizierung', 'SPIDL3', 'SPID Livello 3', 'SPID (National public system of e-ID) Level 3', 'SPID (Nationales öffentliches System der e-ID) Stufe 3'),
  @('2FAHW', 'Autenticazione a due fattori hardware', 'Two Factor Authentication hardware', 'Zwei-Faktor-Hardware-Authentifizierung', 'CNS', 'Carta nazionale dei servizi (CNS)', 'National service card', 'Nationale Servicekarte'),
  @('2FAHW', 'Autenticazione a due fattori hardware', 'Two Factor Authentication hardware', 'Zwei-Faktor-Hardware-Authentifizierung', 'CIE', 'Carta d''identità elettronica (CIE)', 'Electronic identity card', 'Elektronischer Personalausweis (CIE)'),
  @('MFA', 'Autenticazione multi-fattore', 'Multi-factor authentication', 'Multi-Faktor-Authentifizierung', 'USERID OTP BIO', 'User ID, One Time Password e Biometria', 'User ID, One Time Password and Biometric', 'Benutzer-ID, Einmalpasswort und biometrische Daten')
)

for ($r = 0; $r -lt $data.Length; $r++) {
  $rowVals = $data[$r]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $v = $rowVals[$c]
    if ($v -ne $null) {
      $ws.Cells.Item($r + 1, $c + 1).Value = $v
    }
  }
}

# ---------------------------------------------------------------------------
# 2. Re-apply formatting by copying it from the cells that still hold the
#    original style definitions (copy/paste-formats reuses the existing xf
#    index instead of minting a new one).
# ---------------------------------------------------------------------------

# style s="3" (bold header) - already on A1:D1, extend across the new header row
$ws.Range("A1").Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)

# style s="4" - already on A3 (and A4), extend to the new level-2 code column
$ws.Range("A3").Copy()
$u = $excel.Union($ws.Range("A3:A4"), $ws.Range("E3"))
$u.PasteSpecial(-4122)

# style s="5" - already on B2, covers the whole level-1 label block B2:D10
$ws.Range("B2").Copy()
$ws.Range("B2:D10").PasteSpecial(-4122)

# style s="2" - already on D3 (wrap text + bold font), now used by F3
$ws.Range("D3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

# style s="1" - already on D4, now used by E2 (blank/styled) and F4:F10
$ws.Range("D4").Copy()
$u2 = $excel.Union($ws.Range("E2"), $ws.Range("F4:F10"))
$u2.PasteSpecial(-4122)

# style s="6" (new) - default font with wrap text, used by G4:H4
$ws.Range("G4:H4").WrapText = $true

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Column widths / row height for the new columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 21.83203125
$ws.Rows.Item(4).RowHeight = 75

# ---------------------------------------------------------------------------
# 4. View: zoom to 150% and move the selection to H2.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 150
$ws.Range("H2").Select()
